$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "general": update several summary values (P7 -10% test)
# ---------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("general")
$wsGeneral.Range("B3").Value = 56.63098846897222
$wsGeneral.Range("B4").Value = 0.01900005340576172
$wsGeneral.Range("B6").Value = 38.17098846897686
$wsGeneral.Range("B7").Value = 0.6217041096856285
$wsGeneral.Range("B8").Value = 0.6217041096856285
$wsGeneral.Range("B9").Value = 18.45999999999535
$wsGeneral.Range("B10").Value = 0

# ---------------------------------------------------------------
# Sheet "x": column B (j) values change for several rows
# ---------------------------------------------------------------
$wsX = $wb.Worksheets.Item("x")
$wsX.Range("B3").Value = 3
$wsX.Range("B6").Value = 5
$wsX.Range("B7").Value = 4
$wsX.Range("B9").Value = 12
$wsX.Range("B10").Value = 2
$wsX.Range("B13").Value = 13
$wsX.Range("B14").Value = 11

# ---------------------------------------------------------------
# Sheet "U": single cell change
# ---------------------------------------------------------------
$wsU = $wb.Worksheets.Item("U")
$wsU.Range("B2").Value = 3

# ---------------------------------------------------------------
# Sheet "TBar": column B values change for several rows
# ---------------------------------------------------------------
$wsTBar = $wb.Worksheets.Item("TBar")
$wsTBar.Range("B3").Value = 30
$wsTBar.Range("B4").Value = 34.3488504228129
$wsTBar.Range("B5").Value = 30
$wsTBar.Range("B6").Value = 35.00919155153804
$wsTBar.Range("B7").Value = 37.06506101847738
$wsTBar.Range("B9").Value = 30.60033324079214
$wsTBar.Range("B10").Value = 37.32144153802307
$wsTBar.Range("B11").Value = 39.77511225757775
$wsTBar.Range("B12").Value = 32.61192465059682
$wsTBar.Range("B13").Value = 36.71671453559702
$wsTBar.Range("B14").Value = 38.25017704655227
$wsTBar.Range("B15").Value = 42.88363280600589

# ---------------------------------------------------------------
# Sheet "y": rows 3-9 removed, row 2 updated (A1:D9 -> A1:D2)
# ---------------------------------------------------------------
$wsY = $wb.Worksheets.Item("y")
$wsY.Range("A3:D9").EntireRow.Delete()
$wsY.Range("A2").Value = 9
$wsY.Range("B2").Value = 13
$wsY.Range("C2").Value = 2
$wsY.Range("D2").Value = 1

# ---------------------------------------------------------------
# Sheet "Q": column C values change for rows 7-71
# ---------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("Q")
$wsQ.Range("C7").Value = 109.9450000000008
$wsQ.Range("C8").Value = 117.5900000000008
$wsQ.Range("C9").Value = 113.2700000000008
$wsQ.Range("C10").Value = 119.1550000000008
$wsQ.Range("C11").Value = 115.8050000000008
$wsQ.Range("C12").Value = 188.8550000000006
$wsQ.Range("C13").Value = 192.9200000000006
$wsQ.Range("C14").Value = 178.5050000000006
$wsQ.Range("C15").Value = 189.2700000000006
$wsQ.Range("C16").Value = 182.1250000000006
$wsQ.Range("C17").Value = 46.91999999999942
$wsQ.Range("C18").Value = 36.10499999999942
$wsQ.Range("C19").Value = 34.91499999999942
$wsQ.Range("C20").Value = 37.48999999999942
$wsQ.Range("C21").Value = 39.43499999999941
$wsQ.Range("C22").Value = 119.55
$wsQ.Range("C23").Value = 116.1599999999989
$wsQ.Range("C24").Value = 117.2349999999989
$wsQ.Range("C25").Value = 121.4449999999989
$wsQ.Range("C26").Value = 120.25
$wsQ.Range("C27").Value = 258.7350000000008
$wsQ.Range("C28").Value = 269.2400000000008
$wsQ.Range("C29").Value = 250.9150000000008
$wsQ.Range("C30").Value = 261.9150000000008
$wsQ.Range("C31").Value = 255.0150000000008
$wsQ.Range("C37").Value = 193.0200000000017
$wsQ.Range("C38").Value = 202.3100000000017
$wsQ.Range("C39").Value = 191.2450000000017
$wsQ.Range("C40").Value = 208.9250000000017
$wsQ.Range("C41").Value = 197.6600000000017
$wsQ.Range("C42").Value = 187.4749999999983
$wsQ.Range("C43").Value = 195.3199999999983
$wsQ.Range("C44").Value = 177.0549999999983
$wsQ.Range("C45").Value = 185.2149999999984
$wsQ.Range("C46").Value = 179.1799999999983
$wsQ.Range("C47").Value = 272.9599999999988
$wsQ.Range("C48").Value = 283.2849999999988
$wsQ.Range("C49").Value = 256.7699999999988
$wsQ.Range("C50").Value = 275.9449999999989
$wsQ.Range("C51").Value = 263.9099999999988
$wsQ.Range("C52").Value = 250.970000000001
$wsQ.Range("C53").Value = 260.9900000000009
$wsQ.Range("C54").Value = 252.975000000001
$wsQ.Range("C55").Value = 269.580000000001
$wsQ.Range("C56").Value = 250.575000000001
$wsQ.Range("C57").Value = 250.970000000001
$wsQ.Range("C58").Value = 260.9900000000009
$wsQ.Range("C59").Value = 252.975000000001
$wsQ.Range("C60").Value = 269.580000000001
$wsQ.Range("C61").Value = 250.575000000001
$wsQ.Range("C62").Value = 258.7350000000008
$wsQ.Range("C63").Value = 269.2400000000008
$wsQ.Range("C64").Value = 250.9150000000008
$wsQ.Range("C65").Value = 261.9150000000008
$wsQ.Range("C66").Value = 255.0150000000008
$wsQ.Range("C67").Value = 272.9599999999988
$wsQ.Range("C68").Value = 283.2849999999988
$wsQ.Range("C69").Value = 256.7699999999988
$wsQ.Range("C70").Value = 275.9449999999989
$wsQ.Range("C71").Value = 263.9099999999988

# ---------------------------------------------------------------
# Sheet "R": column C values change for rows 7-16
# ---------------------------------------------------------------
$wsR = $wb.Worksheets.Item("R")
$wsR.Range("C7").Value = 0
$wsR.Range("C8").Value = 0
$wsR.Range("C9").Value = 0
$wsR.Range("C10").Value = 0
$wsR.Range("C11").Value = 0
$wsR.Range("C12").Value = 0
$wsR.Range("C13").Value = 8.284999999998837
$wsR.Range("C14").Value = 0
$wsR.Range("C15").Value = 0.9449999999988359
$wsR.Range("C16").Value = 0

# ---------------------------------------------------------------
# Sheet "L": column C values change for rows 2-6
# ---------------------------------------------------------------
$wsL = $wb.Worksheets.Item("L")
$wsL.Range("C2").Value = 0
$wsL.Range("C3").Value = 0
$wsL.Range("C4").Value = 0
$wsL.Range("C5").Value = 0
$wsL.Range("C6").Value = 0

# ---------------------------------------------------------------
# Sheet "rho": rows 4-11 removed, rows 2-3 updated (A1:C11 -> A1:C3)
# ---------------------------------------------------------------
$wsRho = $wb.Worksheets.Item("rho")
$wsRho.Range("A4:C11").EntireRow.Delete()
$wsRho.Range("A2").Value = 9
$wsRho.Range("B2").Value = 2
$wsRho.Range("C2").Value = 1
$wsRho.Range("A3").Value = 9
$wsRho.Range("B3").Value = 4
$wsRho.Range("C3").Value = 1

# ---------------------------------------------------------------
# Sheet "alpha": rows 3-9 removed, row 2 updated (A1:C9 -> A1:C2)
# ---------------------------------------------------------------
$wsAlpha = $wb.Worksheets.Item("alpha")
$wsAlpha.Range("A3:C9").EntireRow.Delete()
$wsAlpha.Range("A2").Value = 9
$wsAlpha.Range("B2").Value = 2
$wsAlpha.Range("C2").Value = 1
